$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: RandomForestRegressor (name unchanged), update metric values
$ws.Range("B3").Value = 0.990422046723012
$ws.Range("C3").Value = 0.9900990409403017
$ws.Range("D3").Value = 0.990582694946029

# Row 4: rename model and update metric values
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.9908441262156964
$ws.Range("C4").Value = 0.990847269271228
$ws.Range("D4").Value = 0.9912478174339981

# Row 5: rename model and update metric values
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.874680470187448
$ws.Range("C5").Value = 0.8543316757481996
$ws.Range("D5").Value = 0.8536062902627185
